{"js": "// The worksheet holds a single table of two-digit \u00f7 one-digit division\n// prompts (\"NN\u00f7N=\") laid out row-major (only every 4th row actually has\n// text; the rest are blank answer rows). The commit replaces each of the\n// 25 prompts, in document/reading order, with a new prompt while leaving\n// every other part of the document (the date line, formatting, blank\n// rows, ...) untouched. We walk the table's non-empty cells in row-major\n// order and swap in the Nth replacement, which keeps this robust even if\n// the \"which rows hold data\" layout were to shift.\n\nconst replacements = [\n  \"61\u00f77=\", \"14\u00f78=\", \"38\u00f73=\", \"19\u00f72=\", \"93\u00f74=\",\n  \"91\u00f76=\", \"98\u00f73=\", \"40\u00f77=\", \"43\u00f76=\", \"61\u00f74=\",\n  \"64\u00f73=\", \"49\u00f76=\", \"49\u00f72=\", \"91\u00f75=\", \"28\u00f76=\",\n  \"71\u00f72=\", \"15\u00f76=\", \"40\u00f74=\", \"97\u00f72=\", \"38\u00f74=\",\n  \"18\u00f73=\", \"83\u00f73=\", \"76\u00f78=\", \"63\u00f76=\", \"22\u00f78=\"\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst values = table.values;\nlet replIndex = 0;\n\nfor (let r = 0; r < values.length && replIndex < replacements.length; r++) {\n  const row = values[r];\n  for (let c = 0; c < row.length && replIndex < replacements.length; c++) {\n    if (row[c] !== \"\") {\n      const cell = table.getCell(r, c);\n      const para = cell.body.paragraphs.getFirst();\n      para.getRange().insertText(replacements[replIndex], Word.InsertLocation.replace);\n      replIndex++;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# The worksheet holds a single table of two-digit \u00f7 one-digit division\n# prompts (\"NN\u00f7N=\") laid out row-major (only every 4th row actually has\n# text; the rest are blank answer rows). The commit replaces each of the\n# 25 prompts, in document/reading order, with a new prompt while leaving\n# every other part of the document (the date line, formatting, blank\n# rows, ...) untouched. We walk the table's non-empty cells in row-major\n# order and swap in the Nth replacement, which keeps this robust even if\n# the \"which rows hold data\" layout were to shift.\n\n$replacements = @(\n    \"61\u00f77=\", \"14\u00f78=\", \"38\u00f73=\", \"19\u00f72=\", \"93\u00f74=\",\n    \"91\u00f76=\", \"98\u00f73=\", \"40\u00f77=\", \"43\u00f76=\", \"61\u00f74=\",\n    \"64\u00f73=\", \"49\u00f76=\", \"49\u00f72=\", \"91\u00f75=\", \"28\u00f76=\",\n    \"71\u00f72=\", \"15\u00f76=\", \"40\u00f74=\", \"97\u00f72=\", \"38\u00f74=\",\n    \"18\u00f73=\", \"83\u00f73=\", \"76\u00f78=\", \"63\u00f76=\", \"22\u00f78=\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replIndex = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        if ($replIndex -ge $replacements.Count) { break }\n        $cell = $t.Cell($r, $c)\n        # Cell.Range.Text includes the trailing cell-mark character(s);\n        # strip them so we can tell whether the cell actually has a\n        # division prompt in it before touching it.\n        $cellText = $cell.Range.Text.TrimEnd([char]0x07, [char]0x0D)\n        if ($cellText -ne \"\") {\n            $cell.Range.Text = $replacements[$replIndex]\n            $replIndex++\n        }\n    }\n}\n"}
